$wb = $excel.ActiveWorkbook

# Update the status text on all three sheets: "Ready for handoff" -> "In Translation"
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# Recompute column widths to reflect the shorter status text
# (target stored width is ~13.41 chars; the engine quantizes ColumnWidth to
# 1/6-character steps, so 12.5 is the input that lands closest on that grid)
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
